$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '75.901.02'
$ws.Range("E2").Value = '  +1.13%  '

$ws.Range("D3").Value = '2.900.97'
$ws.Range("E3").Value = '  +1.75%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '198.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '597.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.20%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -1.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.195'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.49%  '

$ws.Range("D10").Value = '2.898.76'
$ws.Range("E10").Value = '  +1.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.419'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +12.29%  '

$ws.Range("E12").Value = '  -1.47%  '

$ws.Range("E13").Value = '  -1.60%  '

$ws.Range("D14").Value = '3.422.58'
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").Value = '75.806.34'
$ws.Range("E15").Value = '  +0.99%  '

$ws.Range("E16").Value = '  -0.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.03%  '

$ws.Range("D18").Value = '2.905.74'
$ws.Range("E18").Value = '  +1.41%  '

$ws.Range("E19").Value = '  -3.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.45%  '

$ws.Range("E23").Value = '  -0.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.31%  '

$ws.Range("E26").Value = '  +1.18%  '

$ws.Range("E27").Value = '  -1.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.69%  '

$ws.Range("E29").Value = '  +2.90%  '

$ws.Range("E30").Value = '  -1.41%  '

$ws.Range("E31").Value = '  -1.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '502.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.66%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.71'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.45%  '

$ws.Range("E34").Value = '  -2.20%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.08%  '

$ws.Range("E39").Value = '  -7.43%  '

$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '179.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.341'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.14%  '

$ws.Range("E44").Value = '  -2.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0904'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.70%  '

$ws.Range("E46").Value = '  -5.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.57%  '

$ws.Range("E49").Value = '  -0.57%  '

$ws.Range("E50").Value = '  -1.73%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.653'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.65%  '
